# Rename IndividualPhysiology.xlsx -> IndividualBiometrics.xlsx
# (row 6 of Tabelle1: individualPhysiologyFile / value / description)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("B6").Value = "IndividualBiometrics.xlsx"
$ws.Range("C6").Value = "Name of the excel file with individual biometrics information. Must be located in the ""paramsFolder"""

[void]$ws.Range("C6").Select()
